$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: wipe all existing cell contents (values + the old shared
# strings they referenced) so the sheet can be rebuilt from scratch while
# preserving the header row's formatting (style id 1 on A1:G1).
$ws.UsedRange.ClearContents()

# Give the new header cell H1 the same look as the rest of the header row
# (copy style from G1) before we put a value into it.
$ws.Range("G1").Copy($ws.Range("H1"))

# ---- Header row (row 1) -------------------------------------------------
$ws.Range("B1").Value = "Total [n]"
$ws.Range("D1").Value = "Funktion erlebbar [n]"
$ws.Range("E1").Value = "Funktion erlebbar [%]"
$ws.Range("F1").Value = "min"
$ws.Range("G1").Value = "Schnitt"
$ws.Range("H1").Value = "max"

# ---- Label-only rows (column A) -----------------------------------------
$ws.Range("A2").Value = "MIB Resets"
$ws.Range("A4").Value = "Hauptkontexte erreicht"
$ws.Range("A5").Value = "demnach nicht erreichte Hauptkontexte"
$ws.Range("A6").Value = "davon Klangkontext"
$ws.Range("A7").Value = "davon Phonekontext"
$ws.Range("A8").Value = "davon Multimedia kontext"
$ws.Range("A9").Value = "davon Car Kontext"
$ws.Range("A10").Value = "davon App Connect Context"
$ws.Range("A11").Value = "davon Navigations Kontext"
$ws.Range("A12").Value = "davon Setup kontext"
$ws.Range("A14").Value = "VIN im Fahrzeugkontext"
$ws.Range("A15").Value = "Legal Kontext erreichbar"
$ws.Range("A16").Value = "Legal Kontext lädt Info vom Backend"
$ws.Range("A17").Value = "Shop Kontext erreichbar"
$ws.Range("A18").Value = "Apps im Shop verfügbar"
$ws.Range("A21").Value = "Audio nach Aufstart hoerbar (USB-Stick)"
$ws.Range("A22").Value = "Bild Player sichtbar (USB-Stick)"
$ws.Range("A24").Value = "Audio nach Aufstart hoerbar (USB-iOS)"
$ws.Range("A25").Value = "Bild Player sichtbar (USB-iOS)"
$ws.Range("A27").Value = "Audio nach Aufstart hoerbar (USB-Android)"
$ws.Range("A28").Value = "Bild Player sichtbar (USB-Android)"

# ---- Result rows with calculated values, mapped onto the new headers ----
$ws.Range("A30").Value = "First result info"
$ws.Range("B30").Value = 9
$ws.Range("D30").Value = 89
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = 15

$ws.Range("A32").Value = "Second result info"
$ws.Range("D32").Value = 100
$ws.Range("F32").Value = 4
$ws.Range("G32").Value = 4
$ws.Range("H32").Value = 4

$ws.Range("A34").Value = "Third result info"
$ws.Range("D34").Value = 100

$ws.Range("A35").Value = "Navigation"
$ws.Range("B35").Value = 1
